$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.915.90"
$ws.Range("E2").Value = "  +1.08%  "
$ws.Range("D3").Value = "1.709.13"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.23"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4028"
$ws.Range("E7").Value = "  +3.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4067"
$ws.Range("E8").Value = "  +0.60%  "
$ws.Range("B9").Value = "BinanceUSD"
$ws.Range("C9").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.006"
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.481"
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.84"
$ws.Range("E11").Value = "  +1.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08823"
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.48"
$ws.Range("E13").Value = "  +7.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.502"
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.034"
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001343"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").Value = "1.699.07"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.05"
$ws.Range("E18").Value = "  -3.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07179"
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.93"
$ws.Range("E20").Value = "  +5.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.272"
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.004"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("E23").Value = "  +1.76%  "
$ws.Range("D24").Value = "24.909.56"
$ws.Range("E24").Value = "  +1.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.337"
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.892"
$ws.Range("E26").Value = "  -4.04%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.14"
$ws.Range("E27").Value = "  +1.73%  "
$ws.Range("B28").Value = "HuobiToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.412"
$ws.Range("E28").Value = "  +22.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.31"
$ws.Range("E29").Value = "  +0.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "143.67"
$ws.Range("E30").Value = "  +5.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.195"
$ws.Range("E31").Value = "  -2.41%  "
$ws.Range("E32").Value = "  +13.88%  "
$ws.Range("D33").Value = "1.895.23"
$ws.Range("E33").Value = "  +1.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08754"
$ws.Range("E34").Value = "  -1.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.03185"
$ws.Range("E35").Value = "  +8.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.329"
$ws.Range("E36").Value = "  -2.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.032"
$ws.Range("E37").Value = "  -0.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2879"
$ws.Range("E38").Value = "  +5.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8468"
$ws.Range("E39").Value = "  +7.88%  "
$ws.Range("E40").Value = "  +1.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09468"
$ws.Range("E41").Value = "  +3.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.14"
$ws.Range("E42").Value = "  -1.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.475"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.69"
$ws.Range("E44").Value = "  +5.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.726"
$ws.Range("E45").Value = "  +5.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7460"
$ws.Range("E46").Value = "  +3.46%  "
$ws.Range("E47").Value = "  +0.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.390"
$ws.Range("E48").Value = "  +4.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.003"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.25"
$ws.Range("E50").Value = "  +1.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.08423"
$ws.Range("E51").Value = "  +5.70%  "
